# Fruta / hortaliza, semanal
#
# A new weekly record for "Chirimoya" (Vega Modelo de Temuco) needs to be
# inserted as row 19 (Fecha 2021-09-08 / serial 44447), pushing the
# existing rows 19-42 down to 20-43 (dimension grows from T42 to T43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 19; this shifts rows 19:42 down
# to 20:43 and carries the column-D date style (s="2") along with it,
# exactly like Excel's native "Insert Sheet Rows" command.
$ws.Rows("19:19").Insert()

# Populate the newly-inserted row 19 with the new observation.
$ws.Cells.Item(19, 1).Value  = 10
$ws.Cells.Item(19, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value  = "La Araucanía"
$ws.Cells.Item(19, 4).Value  = 44447
$ws.Cells.Item(19, 5).Value  = 9
$ws.Cells.Item(19, 6).Value  = "Fruta"
$ws.Cells.Item(19, 7).Value  = 100107
$ws.Cells.Item(19, 8).Value  = "Otros"
$ws.Cells.Item(19, 9).Value  = 100107002
$ws.Cells.Item(19, 10).Value = "Chirimoya"
$ws.Cells.Item(19, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 40
$ws.Cells.Item(19, 14).Value = 3000
$ws.Cells.Item(19, 15).Value = 3000
$ws.Cells.Item(19, 16).Value = 3000
$ws.Cells.Item(19, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(19, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(19, 19).Value = 3000
$ws.Cells.Item(19, 20).Value = 1
